$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rescale the data columns.
#    Columns B,C,D,E (2-5) are divided by 64.
#    Columns G,H,I,J,K,L,M,N (7-14) are divided by 8.
#    Columns A (rank) and F (percentage) are left untouched.
# ---------------------------------------------------------------------------
$div64Cols = @(2, 3, 4, 5)
$div8Cols  = @(7, 8, 9, 10, 11, 12, 13, 14)

for ($r = 2; $r -le 37; $r++) {
    foreach ($c in $div64Cols) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $cell.Value2 / 64
    }
    foreach ($c in $div8Cols) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $cell.Value2 / 8
    }
}

# ---------------------------------------------------------------------------
# 2. Clear the stray explicit (font-only) formatting on C6 so it goes back to
#    the default style, matching the cleaned-up style table.
# ---------------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Update the sheet view: scrolled so column D is the left-most visible
#    column, and the active selection is now P3.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("P3").Select()

# ---------------------------------------------------------------------------
# 4. Widen the columns that now show the rescaled (decimal) values so they
#    display without truncation.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 10.85546875
$ws.Columns.Item(7).ColumnWidth = 12.140625
$ws.Columns.Item(9).ColumnWidth = 11.5703125
$ws.Columns.Item(10).ColumnWidth = 14
$ws.Columns.Item(11).ColumnWidth = 15.28515625
$ws.Columns.Item(12).ColumnWidth = 13.5703125
$ws.Columns.Item(13).ColumnWidth = 11.5703125
